# Planilha-Produtos: rework headers/content on both sheets + restyle the
# header row (smaller font, darker fill).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Produtos")
$ws2 = $wb.Worksheets.Item("Ofertas")

# ---------------------------------------------------------------------
# Produtos sheet: new header row (A:D) + two simple product name rows
# ---------------------------------------------------------------------
$ws1.Range("A1").Value = "Nome do Produto"
$ws1.Range("B1").Value = "Quantidade"
$ws1.Range("C1").Value = "Vendidos"
$ws1.Range("D1").Value = "Preço"

$ws1.Range("A2").Value = "Bacon T"
$ws1.Range("B2").ClearContents()
$ws1.Range("C2").ClearContents()

$ws1.Range("A3").Value = "Cama X"

# New column D needs the same width + header styling as A:C
$ws1.Columns.Item(4).ColumnWidth = 29.2
$ws1.Range("C1").Copy()
$ws1.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Ofertas sheet: trimmed down to 3 columns, header styled like Produtos,
# plus one (currently empty) data row below the header
# ---------------------------------------------------------------------
$ws2.Range("A1").Value = "Nome da Oferta"
$ws2.Range("B1").Value = "Data Inicio"
$ws2.Range("C1").Value = "Data Fim"
$ws2.Range("D1").ClearContents()

$ws2.Columns.Item(1).ColumnWidth = 29.2
$ws2.Columns.Item(2).ColumnWidth = 29.2
$ws2.Columns.Item(3).ColumnWidth = 29.2

$ws1.Range("A1:C1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# touch A2 so the (empty) second row exists, without giving it a style
$ws2.Range("A2").WrapText = $false

# ---------------------------------------------------------------------
# Header style rework: smaller font + darker fill (applies to the style
# shared by both header rows)
# ---------------------------------------------------------------------
$ws1.Range("A1:D1").Font.Size = 11
$ws1.Range("A1:D1").Interior.Color = 3355443

$ws2.Range("A1:C1").Font.Size = 11
$ws2.Range("A1:C1").Interior.Color = 3355443

# ---------------------------------------------------------------------
# Selection / active sheet bookkeeping to match the final document state
# (Produtos stays the active/selected tab, each sheet keeps its own
# cursor position)
# ---------------------------------------------------------------------
$ws2.Range("C1").Select()
$ws1.Range("D1").Select()
$ws1.Activate()
